$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.63904266666666
$ws.Range("H2").Value = 112.917128
$ws.Range("I2").Value = 0.4850220755088102
$ws.Range("J2").Value = 0.4850220755088102
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.752878
$ws.Range("N2").Value = 5.258634
$ws.Range("O2").Value = 0.1377607590022273
$ws.Range("P2").Value = 0.1377607590022273
$ws.Range("Q2").Value = 65.97664983146132
$ws.Range("R2").Value = 593.7898484831519
$ws.Range("S2").Value = 0.06681700925492928
$ws.Range("T2").Value = 0.06681700925492928
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.63904266666666
$ws.Range("H3").Value = 112.917128
$ws.Range("I3").Value = 0.4850220755088102
$ws.Range("J3").Value = 0.4850220755088102
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.076282333333333
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.2417686736584878
$ws.Range("P3").Value = 0.2417686736584878
$ws.Range("Q3").Value = 115.7883219990462
$ws.Range("R3").Value = 1042.094897991416
$ws.Range("S3").Value = 0.1172631438908519
$ws.Range("T3").Value = 0.117263143890852
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.63904266666666
$ws.Range("H4").Value = 112.917128
$ws.Range("I4").Value = 0.4850220755088102
$ws.Range("J4").Value = 0.4850220755088102
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.303088666666667
$ws.Range("N4").Value = 21.909266
$ws.Range("O4").Value = 0.5739583917309499
$ws.Range("P4").Value = 0.5739583917309499
$ws.Range("Q4").Value = 274.8812659231164
$ws.Range("R4").Value = 2473.931393308048
$ws.Range("S4").Value = 0.278382490413044
$ws.Range("T4").Value = 0.278382490413044
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.63904266666666
$ws.Range("H5").Value = 112.917128
$ws.Range("I5").Value = 0.4850220755088102
$ws.Range("J5").Value = 0.4850220755088102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5918243333333334
$ws.Range("N5").Value = 1.775473
$ws.Range("O5").Value = 0.04651217560833507
$ws.Range("P5").Value = 0.04651217560833507
$ws.Range("Q5").Value = 22.27570133350489
$ws.Range("R5").Value = 200.481312001544
$ws.Range("S5").Value = 0.02255943194998493
$ws.Range("T5").Value = 0.02255943194998494
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.57434666666667
$ws.Range("H6").Value = 52.72304
$ws.Range("I6").Value = 0.2264655392929762
$ws.Range("J6").Value = 0.2264655392929762
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.752878
$ws.Range("N6").Value = 5.258634
$ws.Range("O6").Value = 0.1377607590022273
$ws.Range("P6").Value = 0.1377607590022273
$ws.Range("Q6").Value = 30.80568563637333
$ws.Range("R6").Value = 277.25117072736
$ws.Range("S6").Value = 0.03119806458084913
$ws.Range("T6").Value = 0.03119806458084913
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.57434666666667
$ws.Range("H7").Value = 52.72304
$ws.Range("I7").Value = 0.2264655392929762
$ws.Range("J7").Value = 0.2264655392929762
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.076282333333333
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.2417686736584878
$ws.Range("P7").Value = 0.2417686736584878
$ws.Range("Q7").Value = 54.06365217054223
$ws.Range("R7").Value = 486.57286953488
$ws.Range("S7").Value = 0.054752273064217
$ws.Range("T7").Value = 0.05475227306421701
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.57434666666667
$ws.Range("H8").Value = 52.72304
$ws.Range("I8").Value = 0.2264655392929762
$ws.Range("J8").Value = 0.2264655392929762
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.303088666666667
$ws.Range("N8").Value = 21.909266
$ws.Range("O8").Value = 0.5739583917309499
$ws.Range("P8").Value = 0.5739583917309499
$ws.Range("Q8").Value = 128.3470119654044
$ws.Range("R8").Value = 1155.12310768864
$ws.Range("S8").Value = 0.1299817967150788
$ws.Range("T8").Value = 0.1299817967150788
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.57434666666667
$ws.Range("H9").Value = 52.72304
$ws.Range("I9").Value = 0.2264655392929762
$ws.Range("J9").Value = 0.2264655392929762
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5918243333333334
$ws.Range("N9").Value = 1.775473
$ws.Range("O9").Value = 0.04651217560833507
$ws.Range("P9").Value = 0.04651217560833507
$ws.Range("Q9").Value = 10.40092599976889
$ws.Range("R9").Value = 93.60833399792
$ws.Range("S9").Value = 0.01053340493283121
$ws.Range("T9").Value = 0.01053340493283121
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.168158
$ws.Range("H10").Value = 57.504474
$ws.Range("I10").Value = 0.2470036195972184
$ws.Range("J10").Value = 0.2470036195972184
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.752878
$ws.Range("N10").Value = 5.258634
$ws.Range("O10").Value = 0.1377607590022273
$ws.Range("P10").Value = 0.1377607590022273
$ws.Range("Q10").Value = 33.599442458724
$ws.Range("R10").Value = 302.394982128516
$ws.Range("S10").Value = 0.03402740611201023
$ws.Range("T10").Value = 0.03402740611201023
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 19.168158
$ws.Range("H11").Value = 57.504474
$ws.Range("I11").Value = 0.2470036195972184
$ws.Range("J11").Value = 0.2470036195972184
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.076282333333333
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.2417686736584878
$ws.Range("P11").Value = 0.2417686736584878
$ws.Range("Q11").Value = 58.96666581794201
$ws.Range("R11").Value = 530.6999923614781
$ws.Range("S11").Value = 0.05971773749886515
$ws.Range("T11").Value = 0.05971773749886516
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 19.168158
$ws.Range("H12").Value = 57.504474
$ws.Range("I12").Value = 0.2470036195972184
$ws.Range("J12").Value = 0.2470036195972184
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.303088666666667
$ws.Range("N12").Value = 21.909266
$ws.Range("O12").Value = 0.5739583917309499
$ws.Range("P12").Value = 0.5739583917309499
$ws.Range("Q12").Value = 139.986757450676
$ws.Range("R12").Value = 1259.880817056084
$ws.Range("S12").Value = 0.1417698002557428
$ws.Range("T12").Value = 0.1417698002557428
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 19.168158
$ws.Range("H13").Value = 57.504474
$ws.Range("I13").Value = 0.2470036195972184
$ws.Range("J13").Value = 0.2470036195972184
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.5918243333333334
$ws.Range("N13").Value = 1.775473
$ws.Range("O13").Value = 0.04651217560833507
$ws.Range("P13").Value = 0.04651217560833507
$ws.Range("Q13").Value = 11.344182329578
$ws.Range("R13").Value = 102.097640966202
$ws.Range("S13").Value = 0.01148867573060022
$ws.Range("T13").Value = 0.01148867573060022
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.221194
$ws.Range("H14").Value = 9.663582
$ws.Range("I14").Value = 0.04150876560099527
$ws.Range("J14").Value = 0.04150876560099527
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.752878
$ws.Range("N14").Value = 5.258634
$ws.Range("O14").Value = 0.1377607590022273
$ws.Range("P14").Value = 0.1377607590022273
$ws.Range("Q14").Value = 5.646360096332
$ws.Range("R14").Value = 50.81724086698799
$ws.Range("S14").Value = 0.005718279054438651
$ws.Range("T14").Value = 0.005718279054438651
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.221194
$ws.Range("H15").Value = 9.663582
$ws.Range("I15").Value = 0.04150876560099527
$ws.Range("J15").Value = 0.04150876560099527
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.076282333333333
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.2417686736584878
$ws.Range("P15").Value = 0.2417686736584878
$ws.Range("Q15").Value = 9.909302194439334
$ws.Range("R15").Value = 89.18371974995399
$ws.Range("S15").Value = 0.01003551920455369
$ws.Range("T15").Value = 0.01003551920455369
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.221194
$ws.Range("H16").Value = 9.663582
$ws.Range("I16").Value = 0.04150876560099527
$ws.Range("J16").Value = 0.04150876560099527
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.303088666666667
$ws.Range("N16").Value = 21.909266
$ws.Range("O16").Value = 0.5739583917309499
$ws.Range("P16").Value = 0.5739583917309499
$ws.Range("Q16").Value = 23.52466539453467
$ws.Range("R16").Value = 211.721988550812
$ws.Range("S16").Value = 0.02382430434708422
$ws.Range("T16").Value = 0.02382430434708422
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.221194
$ws.Range("H17").Value = 9.663582
$ws.Range("I17").Value = 0.04150876560099527
$ws.Range("J17").Value = 0.04150876560099527
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.5918243333333334
$ws.Range("N17").Value = 1.775473
$ws.Range("O17").Value = 0.04651217560833507
$ws.Range("P17").Value = 0.04651217560833507
$ws.Range("Q17").Value = 1.906380991587334
$ws.Range("R17").Value = 17.157428924286
$ws.Range("S17").Value = 0.00193066299491871
$ws.Range("T17").Value = 0.00193066299491871
